$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds the last-changed date for every data row
# (rows 2-347). Update the whole column from 45175 (2023-09-06) to
# 45177 (2023-09-08), keeping the existing date formatting/style intact.
$ws.Range("C2:C347").Value = 45177
